$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2000
$ws.Cells.Item(40, 9).Value = 2000
$ws.Cells.Item(40, 11).Value = 2000
$ws.Cells.Item(40, 13).Value = -1825
$ws.Cells.Item(58, 8).Value = 2023.3334
$ws.Cells.Item(58, 10).Value = 6000
$ws.Cells.Item(58, 12).Value = 18000
$ws.Cells.Item(58, 14).Value = -18300
$ws.Cells.Item(70, 8).Value = 3759.5334
$ws.Cells.Item(70, 9).Value = 1087.75
$ws.Cells.Item(70, 10).Value = 4731.091
$ws.Cells.Item(70, 11).Value = 3263.25
$ws.Cells.Item(70, 12).Value = 14193.273
$ws.Cells.Item(70, 13).Value = -2993.25
$ws.Cells.Item(70, 14).Value = -14733.273
$ws.Cells.Item(73, 8).Value = 3759.5334
$ws.Cells.Item(73, 9).Value = 1087.75
$ws.Cells.Item(73, 10).Value = 4731.091
$ws.Cells.Item(73, 11).Value = 3263.25
$ws.Cells.Item(73, 12).Value = 14193.273
$ws.Cells.Item(73, 13).Value = -2327.25
$ws.Cells.Item(73, 14).Value = -16065.273
$ws.Cells.Item(87, 8).Value = 89999
$ws.Cells.Item(87, 10).Value = 89999
$ws.Cells.Item(87, 12).Value = 89999
$ws.Cells.Item(87, 14).Value = -92495
$ws.Cells.Item(90, 8).Value = 89999
$ws.Cells.Item(90, 10).Value = 89999
$ws.Cells.Item(90, 12).Value = 269997
$ws.Cells.Item(90, 14).Value = -282477
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1702
$ws.Cells.Item(2, 9).Value = 1337
$ws.Cells.Item(2, 10).Value = 2249.5
$ws.Cells.Item(2, 11).Value = 1337
$ws.Cells.Item(2, 12).Value = 2249.5
$ws.Cells.Item(2, 13).Value = -1224
$ws.Cells.Item(2, 14).Value = -2475.5
$ws.Cells.Item(32, 8).Value = 24947.516
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 24947.516
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 24947.516
$ws.Cells.Item(32, 13).Value = $null
$ws.Cells.Item(32, 14).Value = -25521.516
$ws.Cells.Item(74, 8).Value = 4849.6665
$ws.Cells.Item(74, 9).Value = 2011.5
$ws.Cells.Item(74, 11).Value = 2011.5
$ws.Cells.Item(74, 13).Value = -1137.5
$ws.Cells.Item(77, 8).Value = 4849.6665
$ws.Cells.Item(77, 9).Value = 2011.5
$ws.Cells.Item(77, 11).Value = 10057.5
$ws.Cells.Item(77, 13).Value = -5689.5
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).Value = $null
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).Value = $null
$ws.Cells.Item(116, 8).Value = 1702
$ws.Cells.Item(116, 9).Value = 1337
$ws.Cells.Item(116, 10).Value = 2249.5
$ws.Cells.Item(116, 11).Value = 1337
$ws.Cells.Item(116, 12).Value = 2249.5
$ws.Cells.Item(116, 13).Value = 957
$ws.Cells.Item(116, 14).Value = -6837.5
$ws.Cells.Item(122, 8).Value = 2461.4333
$ws.Cells.Item(122, 9).Value = 1776.7826
$ws.Cells.Item(122, 11).Value = 5330.3478
$ws.Cells.Item(122, 13).Value = -2880.3478
$ws.Cells.Item(132, 8).Value = 3116.88
$ws.Cells.Item(132, 9).Value = 2518.4783
$ws.Cells.Item(132, 11).Value = 7555.4349
$ws.Cells.Item(132, 13).Value = -5025.4349
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1702
$ws.Cells.Item(3, 9).Value = 1337
$ws.Cells.Item(3, 10).Value = 2249.5
$ws.Cells.Item(3, 11).Value = 1337
$ws.Cells.Item(3, 12).Value = 2249.5
$ws.Cells.Item(3, 13).Value = -1223
$ws.Cells.Item(3, 14).Value = -2477.5
$ws.Cells.Item(86, 8).Value = 3829.25
$ws.Cells.Item(86, 9).Value = 2006
$ws.Cells.Item(86, 11).Value = 2006
$ws.Cells.Item(86, 13).Value = -883
$ws.Cells.Item(89, 8).Value = 3829.25
$ws.Cells.Item(89, 9).Value = 2006
$ws.Cells.Item(89, 11).Value = 10030
$ws.Cells.Item(89, 13).Value = -4414
$ws.Cells.Item(99, 8).Value = 1189.8572
$ws.Cells.Item(99, 9).Value = 972.3333
$ws.Cells.Item(99, 10).Value = 2495
$ws.Cells.Item(99, 11).Value = 972.3333
$ws.Cells.Item(99, 12).Value = 2495
$ws.Cells.Item(99, 13).Value = 525.6667
$ws.Cells.Item(99, 14).Value = -5491
$ws.Cells.Item(105, 8).Value = 4291.75
$ws.Cells.Item(105, 9).Value = 3808.625
$ws.Cells.Item(105, 10).Value = 4774.875
$ws.Cells.Item(105, 11).Value = 3808.625
$ws.Cells.Item(105, 12).Value = 4774.875
$ws.Cells.Item(105, 13).Value = -2061.625
$ws.Cells.Item(105, 14).Value = -8268.875
$ws.Cells.Item(106, 8).Value = 20623.334
$ws.Cells.Item(106, 10).Value = 20623.334
$ws.Cells.Item(106, 12).Value = 20623.334
$ws.Cells.Item(106, 14).Value = -23147.334
$ws.Cells.Item(134, 8).Value = 4400.4287
$ws.Cells.Item(134, 9).Value = 1906
$ws.Cells.Item(134, 10).Value = 5398.2
$ws.Cells.Item(134, 11).Value = 5718
$ws.Cells.Item(134, 12).Value = 16194.6
$ws.Cells.Item(134, 13).Value = -3183
$ws.Cells.Item(134, 14).Value = -21264.6
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 439.8
$ws.Cells.Item(22, 9).Value = 466.66666
$ws.Cells.Item(22, 10).Value = 399.5
$ws.Cells.Item(22, 11).Value = 466.66666
$ws.Cells.Item(22, 12).Value = 399.5
$ws.Cells.Item(22, 13).Value = -116.66666
$ws.Cells.Item(22, 14).Value = -1099.5
$ws.Cells.Item(86, 8).Value = 16547.545
$ws.Cells.Item(86, 9).Value = 13999
$ws.Cells.Item(86, 10).Value = 16802.4
$ws.Cells.Item(86, 11).Value = 13999
$ws.Cells.Item(86, 12).Value = 16802.4
$ws.Cells.Item(86, 13).Value = -12876
$ws.Cells.Item(86, 14).Value = -19048.4
$ws.Cells.Item(89, 8).Value = 16547.545
$ws.Cells.Item(89, 9).Value = 13999
$ws.Cells.Item(89, 10).Value = 16802.4
$ws.Cells.Item(89, 11).Value = 69995
$ws.Cells.Item(89, 12).Value = 84012
$ws.Cells.Item(89, 13).Value = -64379
$ws.Cells.Item(89, 14).Value = -95244
$ws.Cells.Item(105, 8).Value = 2140.7693
$ws.Cells.Item(105, 9).Value = 1712.125
$ws.Cells.Item(105, 11).Value = 1712.125
$ws.Cells.Item(105, 13).Value = 34.875
$ws.Cells.Item(132, 8).Value = 4092.8333
$ws.Cells.Item(132, 9).Value = 3386.25
$ws.Cells.Item(132, 10).Value = 5506
$ws.Cells.Item(132, 11).Value = 10158.75
$ws.Cells.Item(132, 12).Value = 16518
$ws.Cells.Item(132, 13).Value = -7628.75
$ws.Cells.Item(132, 14).Value = -21578
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(29, 8).Value = 334377.16
$ws.Cells.Item(29, 9).Value = 667533.7
$ws.Cells.Item(29, 10).Value = 1220.6666
$ws.Cells.Item(29, 11).Value = 2002601.1
$ws.Cells.Item(29, 12).Value = 3661.9998
$ws.Cells.Item(29, 13).Value = -2002324.1
$ws.Cells.Item(29, 14).Value = -4215.9998
$ws.Cells.Item(36, 8).Value = 2125.7778
$ws.Cells.Item(36, 9).Value = 516.5
$ws.Cells.Item(36, 11).Value = 1549.5
$ws.Cells.Item(36, 13).Value = -1380.5
$ws.Cells.Item(122, 8).Value = 1124.875
$ws.Cells.Item(122, 10).Value = 1499.5
$ws.Cells.Item(122, 12).Value = 13495.5
$ws.Cells.Item(122, 14).Value = -18395.5
$ws.Cells.Item(132, 8).Value = 7833
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 7833
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 70497
$ws.Cells.Item(132, 13).Value = $null
$ws.Cells.Item(132, 14).Value = -75557
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3415.1765
$ws.Cells.Item(113, 9).Value = 1632.5
$ws.Cells.Item(113, 10).Value = 4999.778
$ws.Cells.Item(113, 11).Value = 1632.5
$ws.Cells.Item(113, 12).Value = 4999.778
$ws.Cells.Item(113, 13).Value = 537.5
$ws.Cells.Item(113, 14).Value = -9339.778
$ws.Cells.Item(136, 8).Value = 72290
$ws.Cells.Item(136, 10).Value = 72290
$ws.Cells.Item(136, 12).Value = 216870
$ws.Cells.Item(136, 14).Value = -221970
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4635.6665
$ws.Cells.Item(7, 9).Value = 4953.5
$ws.Cells.Item(7, 11).Value = 4953.5
$ws.Cells.Item(7, 13).Value = -4841.5
$ws.Cells.Item(24, 8).Value = 59998.668
$ws.Cells.Item(24, 10).Value = 59998.668
$ws.Cells.Item(24, 12).Value = 59998.668
$ws.Cells.Item(24, 14).Value = -60684.668
$ws.Cells.Item(46, 8).Value = 2714.6155
$ws.Cells.Item(46, 9).Value = 2299.0908
$ws.Cells.Item(46, 11).Value = 2299.0908
$ws.Cells.Item(46, 13).Value = -2111.0908
$ws.Cells.Item(126, 8).Value = 4635.6665
$ws.Cells.Item(126, 9).Value = 4953.5
$ws.Cells.Item(126, 11).Value = 14860.5
$ws.Cells.Item(126, 13).Value = -12390.5
$ws.Cells.Item(132, 8).Value = 6091.3335
$ws.Cells.Item(132, 9).Value = 4249.5
$ws.Cells.Item(132, 11).Value = 12748.5
$ws.Cells.Item(132, 13).Value = -10218.5
$ws.Cells.Item(136, 8).Value = 4999.6665
$ws.Cells.Item(136, 9).Value = 4999.6665
$ws.Cells.Item(136, 11).Value = 14998.9995
$ws.Cells.Item(136, 13).Value = -12448.9995
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).Value = $null
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = $null
$ws.Cells.Item(21, 14).Value = $null
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 13).Value = $null
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = $null
$ws.Cells.Item(35, 14).Value = $null
$ws.Cells.Item(95, 8).Value = 43391.8
$ws.Cells.Item(95, 10).Value = 43391.8
$ws.Cells.Item(95, 12).Value = 43391.8
$ws.Cells.Item(95, 14).Value = -48883.8
$ws.Cells.Item(100, 8).Value = 1253.0667
$ws.Cells.Item(100, 9).Value = 1259.6
$ws.Cells.Item(100, 10).Value = 1240
$ws.Cells.Item(100, 11).Value = 2519.2
$ws.Cells.Item(100, 12).Value = 2480
$ws.Cells.Item(100, 13).Value = -1978.2
$ws.Cells.Item(100, 14).Value = -3562
$ws.Cells.Item(123, 8).Value = 80390
$ws.Cells.Item(123, 9).Value = 80390
$ws.Cells.Item(123, 11).Value = 80390
$ws.Cells.Item(123, 13).Value = -75490
$ws.Cells.Item(126, 8).Value = 49878.76
$ws.Cells.Item(126, 10).Value = 2469.3845
$ws.Cells.Item(126, 12).Value = 7408.1535
$ws.Cells.Item(126, 14).Value = -12348.1535
$ws.Cells.Item(136, 8).Value = 128321.25
$ws.Cells.Item(136, 9).Value = 2929.3333
$ws.Cells.Item(136, 11).Value = 8787.999899999999
$ws.Cells.Item(136, 13).Value = -6237.999899999999

Write-Output "applied changes"
